$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9557737112045288
$ws.Range("B1").Value = 1.721540808677673
$ws.Range("C1").Value = 4.67745304107666
$ws.Range("D1").Value = 1.413659930229187
$ws.Range("E1").Value = 1.143731832504272
